$wb = $excel.ActiveWorkbook

# --- Rename header columns on existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# --- Header row ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$ws3.Range("A2").Value = 45074.99999999999
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = -14.35172327418276
$ws3.Range("D2").Value = 14.60467839756018

$ws3.Range("A3").Value = 45081.99999999999
$ws3.Range("B3").Value = 2
$ws3.Range("C3").Value = -11.22031407341315
$ws3.Range("D3").Value = 14.82111924887827

$ws3.Range("A4").Value = 45130.99999999999
$ws3.Range("B4").Value = 9
$ws3.Range("C4").Value = -4.01976397915342
$ws3.Range("D4").Value = 22.40896942856109

$ws3.Range("A5").Value = 45165.99999999999
$ws3.Range("B5").Value = 14
$ws3.Range("C5").Value = -0.2847717918132381
$ws3.Range("D5").Value = 26.66596066383477

$ws3.Range("A6").Value = 45179.99999999999
$ws3.Range("B6").Value = 16
$ws3.Range("C6").Value = 2.007401388271684
$ws3.Range("D6").Value = 29.55653649017326

$ws3.Range("A7").Value = 45186.99999999999
$ws3.Range("B7").Value = 17
$ws3.Range("C7").Value = 1.772124002525842
$ws3.Range("D7").Value = 29.76176845809618

$ws3.Range("A8").Value = 45207.99999999999
$ws3.Range("B8").Value = 20
$ws3.Range("C8").Value = 5.887874348019175
$ws3.Range("D8").Value = 34.75717153701662

$ws3.Range("A9").Value = 45214.99999999999
$ws3.Range("B9").Value = 21
$ws3.Range("C9").Value = 6.652172400091469
$ws3.Range("D9").Value = 34.47589033474143

$ws3.Range("A10").Value = 45221.99999999999
$ws3.Range("B10").Value = 22
$ws3.Range("C10").Value = 6.872600046099412
$ws3.Range("D10").Value = 36.43081543088491

$ws3.Range("A11").Value = 45228.99999999999
$ws3.Range("B11").Value = 23
$ws3.Range("C11").Value = 8.317368146611384
$ws3.Range("D11").Value = 36.950667553

$ws3.Range("A12").Value = 45235.99999999999
$ws3.Range("B12").Value = 24
$ws3.Range("C12").Value = 8.976918248510371
$ws3.Range("D12").Value = 38.76463507397703

$ws3.Range("A13").Value = 45242.99999999999
$ws3.Range("B13").Value = 25
$ws3.Range("C13").Value = 10.76274494955496
$ws3.Range("D13").Value = 38.64745409864273

$ws3.Range("A14").Value = 45249.99999999999
$ws3.Range("B14").Value = 26
$ws3.Range("C14").Value = 11.44659330546717
$ws3.Range("D14").Value = 39.7079432488498

$ws3.Range("A15").Value = 45256.99999999999
$ws3.Range("B15").Value = 27
$ws3.Range("C15").Value = 12.84118810152416
$ws3.Range("D15").Value = 40.28445708704427

$ws3.Range("A16").Value = 45263.99999999999
$ws3.Range("B16").Value = 28
$ws3.Range("C16").Value = 14.26206413259698
$ws3.Range("D16").Value = 42.7557401955756

# --- Match header + date formatting to the other sheets ---
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws3.Range("A2:A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore original active sheet/selection ---
$ws1.Activate()
